# Logs/work_logs.xlsx - "questions 2 & 3"
# Adds a new work-log entry (row 32): a new shared-string activity
# "Tackled questions 2 & 3 from data analysis" with a duration of 1.5 hrs.
# The Total (F1, =SUM(C2:C32)) and sheet dimension/selection update as a
# natural consequence of the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry on row 32 (no date in column A, matching the other
# no-date continuation rows such as row 5, 28, 30, ...).
$ws.Cells.Item(32, 2).Value = "Tackled questions 2 & 3 from data analysis"
$ws.Cells.Item(32, 3).Value = 1.5

# Move the active selection the way the author's session ended up
# (one past the newly-added data, column D).
$null = $ws.Range("D35").Select()
